$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 08:12"

# Row 20
$ws.Range("B20").Value = 2799
$ws.Range("C20").Value = 123
$ws.Range("D20").Value = 170
$ws.Range("E20").Value = 2617
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 12

# Row 23
$ws.Range("A23").Value = "Israel"
$ws.Range("B23").Value = 2495
$ws.Range("C23").Value = 126
$ws.Range("D23").Value = 66
$ws.Range("E23").Value = 2424
$ws.Range("F23").Value = 41
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 5

# Row 24
$ws.Range("A24").Value = "Turquia"
$ws.Range("B24").Value = 2433
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 2348
$ws.Range("F24").Value = 136
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 59

# Row 44
$ws.Range("B44").Value = 674
$ws.Range("C44").Value = 17
$ws.Range("D44").Value = 43
$ws.Range("E44").Value = 618
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 13

# Row 66
$ws.Range("A66").Value = "Lituania"
$ws.Range("B66").Value = 290
$ws.Range("C66").Value = 16
$ws.Range("D66").Value = 1
$ws.Range("E66").Value = 285
$ws.Range("F66").Value = 1
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 4

# Row 67
$ws.Range("A67").Value = "Armenia"
$ws.Range("B67").Value = 290
$ws.Range("C67").Value = 25
$ws.Range("D67").Value = 18
$ws.Range("E67").Value = 272
$ws.Range("F67").Value = 6
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 0

# Row 68
$ws.Range("A68").Value = "Nueva Zelanda"
$ws.Range("B68").Value = 283
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 27
$ws.Range("E68").Value = 256
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0

# Row 70
$ws.Range("B70").Value = 243
$ws.Range("C70").Value = 1
$ws.Range("D70").Value = 4
$ws.Range("E70").Value = 236
$ws.Range("F70").Value = 8
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 3

# Row 104
$ws.Range("A104").Value = "Georgia"
$ws.Range("B104").Value = 77
$ws.Range("C104").Value = 2
$ws.Range("D104").Value = 10
$ws.Range("E104").Value = 67
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0

# Row 105
$ws.Range("A105").Value = "Camerun"
$ws.Range("B105").Value = 75
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 2
$ws.Range("E105").Value = 72
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 1

# Row 110
$ws.Range("B110").Value = 65
$ws.Range("C110").Value = 5
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 65
$ws.Range("F110").Value = 4
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 0

# Row 151
$ws.Range("B151").Value = 11
$ws.Range("C151").Value = 1
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 11
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0

# Row 153
$ws.Range("A153").Value = "Haiti"

# Row 154
$ws.Range("A154").Value = "Surinam"

# Row 158
$ws.Range("A158").Value = "Gabon"
$ws.Range("B158").Value = 7
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 0
$ws.Range("E158").Value = 6
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 1

# Row 159
$ws.Range("A159").Value = "Niger"
$ws.Range("B159").Value = 7
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 0
$ws.Range("E159").Value = 6
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 1

# Row 160
$ws.Range("A160").Value = "Namibia"
$ws.Range("B160").Value = 7
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 2
$ws.Range("E160").Value = 5
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 0

# Row 161
$ws.Range("A161").Value = "Benin"
$ws.Range("B161").Value = 6
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 6
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 0

# Row 164
$ws.Range("A164").Value = "Mozambique"

# Row 165
$ws.Range("A165").Value = "Fiyi"

# Row 166
$ws.Range("A166").Value = "Siria"

# Row 176
$ws.Range("A176").Value = "San Bartolome"

# Row 177
$ws.Range("A177").Value = "Liberia"

# Row 178
$ws.Range("A178").Value = "Laos"

# Row 181
$ws.Range("A181").Value = "Santa Lucia"

# Row 182
$ws.Range("A182").Value = "Antigua y Barbuda"

# Row 183
$ws.Range("A183").Value = "San Martin (Parte Holandesa)"

# Row 184
$ws.Range("A184").Value = "Birmania"

# Row 186
$ws.Range("A186").Value = "Zimbabue"

# Row 187
$ws.Range("A187").Value = "Gambia"
$ws.Range("B187").Value = 3
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 2
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 1

# Row 188
$ws.Range("A188").Value = "Nepal"
$ws.Range("B188").Value = 3
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = 2
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

# Row 189
$ws.Range("A189").Value = "Mauritania"

# Row 190
$ws.Range("A190").Value = "San Cristobal y Nieves"

# Row 191
$ws.Range("A191").Value = "Islas Virgenes Britanicas"

# Row 192
$ws.Range("A192").Value = "Butan"

# Row 196
$ws.Range("A196").Value = "Guinea-Bisau"

# Row 197
$ws.Range("A197").Value = "Montserrat"

# Row 198
$ws.Range("A198").Value = "Libia"

# Row 199
$ws.Range("A199").Value = "Somalia"

# Row 200
$ws.Range("A200").Value = "Papua Nueva Guinea"

# Row 201
$ws.Range("A201").Value = "Timor Oriental"

# Row 202
$ws.Range("A202").Value = "Granada"

# Row 203
$ws.Range("A203").Value = "Islas Turcas y Caicos"

# Row 204
$ws.Range("A204").Value = "San Vicente y las Granadinas"
